$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 381-382. This shifts all existing records
# (rows 381..486) down by 2 rows, which reproduces the "every record
# moves one weekly slot later" pattern from the diff: what used to be
# row N is now row N+2, and the two rows that fall off the bottom of
# the original range (485-486) land on the newly-extended rows 487-488.
$ws.Range("A381:A382").EntireRow.Insert()

# Populate the two new rows (381 = Primera, 382 = Segunda) with the new
# weekly record (2023-07-28), mirroring the other fixed columns from the
# surrounding data for this market/category.
$ws.Cells.Item(381, 1).Value = 1
$ws.Cells.Item(381, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(381, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(381, 4).Value = 45135
$ws.Cells.Item(381, 5).Value = 15
$ws.Cells.Item(381, 6).Value = 100114014
$ws.Cells.Item(381, 7).Value = "Betarraga"
$ws.Cells.Item(381, 8).Value = "Sin especificar"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 1200
$ws.Cells.Item(381, 11).Value = 700
$ws.Cells.Item(381, 12).Value = 800
$ws.Cells.Item(381, 13).Value = 750
$ws.Cells.Item(381, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(381, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(381, 16).Value = 188
$ws.Cells.Item(381, 17).Value = 4
$ws.Cells.Item(381, 18).Value = "Hortaliza"

$ws.Cells.Item(382, 1).Value = 1
$ws.Cells.Item(382, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(382, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(382, 4).Value = 45135
$ws.Cells.Item(382, 5).Value = 15
$ws.Cells.Item(382, 6).Value = 100114014
$ws.Cells.Item(382, 7).Value = "Betarraga"
$ws.Cells.Item(382, 8).Value = "Sin especificar"
$ws.Cells.Item(382, 9).Value = "Segunda"
$ws.Cells.Item(382, 10).Value = 1200
$ws.Cells.Item(382, 11).Value = 700
$ws.Cells.Item(382, 12).Value = 800
$ws.Cells.Item(382, 13).Value = 750
$ws.Cells.Item(382, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 150
$ws.Cells.Item(382, 17).Value = 5
$ws.Cells.Item(382, 18).Value = "Hortaliza"
